# Kazakhstan Premier League - atualizacao de bases das ligas
# The two matches in each pair below were mis-ordered; this swaps all
# data columns (B through AD) between the two rows of each pair while
# leaving column A (the running row index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$rowPairs = @(
    @(99, 100),
    @(119, 120),
    @(129, 130),
    @(136, 137),
    @(148, 149),
    @(150, 151),
    @(156, 157)
)

foreach ($pair in $rowPairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    $rowValues1 = @()
    $rowValues2 = @()

    foreach ($col in $columns) {
        $rowValues1 += , ($ws.Range("$col$row1").Value2())
        $rowValues2 += , ($ws.Range("$col$row2").Value2())
    }

    for ($i = 0; $i -lt $columns.Count; $i++) {
        $col = $columns[$i]
        $ws.Range("$col$row1").Value2 = $rowValues2[$i]
        $ws.Range("$col$row2").Value2 = $rowValues1[$i]
    }
}
